$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge title rows 1-2 into a single merged cell (A1:G2) instead of
# two separate merges (A1:G1, A2:G2).
$ws.Range("A1:G1").UnMerge()
$ws.Range("A2:G2").UnMerge()
$ws.Range("A1:G2").Merge()

# Row heights: row1 already has an explicit row height; make it a
# custom height explicitly. Row 29 grows to 24pt.
$ws.Rows.Item(1).RowHeight = 20.25
$ws.Rows.Item(29).RowHeight = 24

# Update the "grand total" label text (drop the trailing colon).
$ws.Range("A29").Value = "Tổng cộng"

# Update the active selection.
$ws.Range("B12").Select()
